$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values per row (columns B, D, E, F, G, I, J, L, M)
$data = @{
    2 = @(20.14898406622666, 9.103772783116531, 11.58465725634919, 59.28311112824244, 3.799807561985865, 33.49926481469001, 8.455793553752775, 17.7474390357584, 22.19491516287845)
    3 = @(19.91696896012324, 8.787789735066696, 10.99632891769642, 59.11773098006928, 3.806726861283785, 33.65908640749824, 8.226498626224645, 17.61413639402893, 22.07287425272835)
    4 = @(19.77793981204041, 8.59189632741643, 10.62257600829984, 59.03565646412804, 3.811185251101532, 33.76445150680171, 8.081128539723206, 17.53656073116909, 22.00276621772091)
    5 = @(19.72220769253986, 8.511735039951855, 10.46733511127227, 59.0070940741988, 3.813055141024926, 33.80919890771148, 8.020780404891415, 17.5060462750766, 21.97542982131174)
    6 = @(19.71301093874078, 8.498408458140581, 10.44138723012389, 59.00264579514389, 3.813368847574282, 33.81673826271723, 8.010694025498124, 17.50104636198404, 21.9709656583058)
    7 = @(19.77718437093118, 8.590816401919072, 10.62049394513854, 59.03525150855828, 3.811210253886403, 33.76504766757516, 8.080319092426116, 17.53614472570821, 22.00239253152195)
    8 = @(20.06831072099707, 8.995292163529303, 11.38451753608863, 59.22204168870627, 3.802149935842095, 33.55286454894107, 8.377708701307926, 17.70060439783912, 22.15184629211426)
    9 = @(20.66366138618157, 9.767969421470667, 12.94170761164079, 59.74308241851836, 3.786035624901043, 33.194573433328, 8.922802476270842, 18.05587311625532, 22.48226736156673)
    10 = @(21.1122194472147, 10.3164410105269, 14.04903400705951, 60.22019368744353, 3.775186408730531, 32.96713553542423, 9.298001782815417, 18.33522376158716, 22.74642249970837)
    11 = @(21.31789890798592, 10.5605838429943, 14.52526647853484, 60.45758442932026, 3.770461949975112, 32.87156800252384, 9.462831948040556, 18.46591030332662, 22.87091817902524)
    12 = @(21.39595175679062, 10.65217820456376, 14.70166858706045, 60.55038281411777, 3.768702945745825, 32.83652543529168, 9.524380959978568, 18.51588326675707, 22.91865519168595)
    13 = @(21.37913533962686, 10.63249111543547, 14.66385203853512, 60.53026830372915, 3.769080446991579, 32.84402128928639, 9.511164314114136, 18.50509968397329, 22.9083482685589)
    14 = @(21.32431740959687, 10.56813700586099, 14.53985798200242, 60.46516096543761, 3.770316635098171, 32.86866197931935, 9.467913176353221, 18.47001209090207, 22.8748338154356)
    15 = @(21.29075953823999, 10.52860422491552, 14.46339588905622, 60.42565820947424, 3.77107774045368, 32.88390480456717, 9.441306727810586, 18.44858203786578, 22.85438160602085)
    16 = @(21.09880477053103, 10.30037022246543, 14.01735923148753, 60.20508734059523, 3.77549938308445, 32.97354094183929, 9.2871096185787, 18.32675280554567, 22.73837105188995)
    17 = @(20.98141534047983, 10.15891848896422, 13.73669947336092, 60.0749722165514, 3.778265735104799, 33.03055999830432, 9.190995933292697, 18.25291475682493, 22.66829075083551)
    18 = @(20.91405329373957, 10.07705921973924, 13.57268222685081, 60.00205066128411, 3.779876735878838, 33.06409827050709, 9.135164703303607, 18.21078681706385, 22.62839182593248)
    19 = @(20.89127474547511, 10.04926001168367, 13.51670361829569, 59.97769055493888, 3.780425614011015, 33.07558091972426, 9.116167770558986, 18.19658274212692, 22.61495393356644)
    20 = @(20.99389586406276, 10.17402867202632, 13.76684403491303, 60.08862491421206, 3.777969197958515, 33.02441330279088, 9.201284445263918, 18.2607398174844, 22.67570875643645)
    21 = @(21.34041477946191, 10.58706327197217, 14.57638474174765, 60.48420597421145, 3.76995272350765, 32.86139319488115, 9.480640860664957, 18.48030528554079, 22.88466196308875)
    22 = @(21.56782485884973, 10.85197557868233, 15.0825287064136, 60.75965123602337, 3.764888491615592, 32.7615400188303, 9.658142986035438, 18.62661352816326, 23.02466827444046)
    23 = @(21.44638787110157, 10.71107309643008, 14.81448291778657, 60.61110263108225, 3.767575448818685, 32.81421753613124, 9.563879257379353, 18.54828044544739, 22.94963918628403)
    24 = @(20.98825301956761, 10.16719902354935, 13.75322397211078, 60.08244666169225, 3.778103198252965, 33.02718986601108, 9.196634795370674, 18.25720110034989, 22.67235385546218)
    25 = @(20.50037099996322, 9.561870708379365, 12.50943926121774, 59.58553917939218, 3.790219839928816, 33.28525582773582, 8.779643782359665, 17.95641902724421, 22.38901748646011)
}

$cols = @(2, 4, 5, 6, 7, 9, 10, 12, 13)  # B, D, E, F, G, I, J, L, M

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value = [double]$vals[$i]
    }
}